$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "42.691.32"
$ws.Cells.Item(2, 5).Value = "  -0.14%  "
$ws.Cells.Item(3, 4).Value = "2.519.03"
$ws.Cells.Item(3, 5).Value = "  -1.38%  "
$ws.Cells.Item(4, 5).Value = "  +0.08%  "
$ws.Cells.Item(5, 4).Value = "'317.57"
$ws.Cells.Item(5, 5).Value = "  +4.76%  "
$ws.Cells.Item(6, 4).Value = "'95.71"
$ws.Cells.Item(6, 5).Value = "  -2.71%  "
$ws.Cells.Item(7, 4).Value = "'0.586"
$ws.Cells.Item(7, 5).Value = "  +2.01%  "
$ws.Cells.Item(8, 5).Value = "  +0.06%  "
$ws.Cells.Item(9, 4).Value = "'0.539"
$ws.Cells.Item(9, 5).Value = "  -1.39%  "
$ws.Cells.Item(10, 4).Value = "'36.12"
$ws.Cells.Item(10, 5).Value = "  -1.76%  "
$ws.Cells.Item(11, 4).Value = "'0.0814"
$ws.Cells.Item(11, 5).Value = "  +0.67%  "
$ws.Cells.Item(12, 5).Value = "  +1.12%  "
$ws.Cells.Item(13, 5).Value = "  -3.04%  "
$ws.Cells.Item(14, 5).Value = "  -1.13%  "
$ws.Cells.Item(15, 4).Value = "'15.50"
$ws.Cells.Item(15, 5).Value = "  +4.79%  "
$ws.Cells.Item(16, 4).Value = "2.527.71"
$ws.Cells.Item(16, 5).Value = "  -0.06%  "
$ws.Cells.Item(17, 4).Value = "'0.860"
$ws.Cells.Item(17, 5).Value = "  -2.20%  "
$ws.Cells.Item(18, 4).Value = "42.680.74"
$ws.Cells.Item(18, 5).Value = "  -0.28%  "
$ws.Cells.Item(19, 4).Value = "'12.93"
$ws.Cells.Item(19, 5).Value = "  -2.59%  "
$ws.Cells.Item(20, 5).Value = "  -0.95%  "
$ws.Cells.Item(21, 4).Value = "'6.60"
$ws.Cells.Item(21, 5).Value = "  +0.27%  "
$ws.Cells.Item(22, 4).Value = "'71.57"
$ws.Cells.Item(22, 5).Value = "  -0.07%  "
$ws.Cells.Item(23, 4).Value = "'253.30"
$ws.Cells.Item(23, 5).Value = "  -0.50%  "
$ws.Cells.Item(24, 4).Value = "'2.98"
$ws.Cells.Item(24, 5).Value = "  +1.05%  "
$ws.Cells.Item(25, 4).Value = "'2.04"
$ws.Cells.Item(25, 5).Value = "  -2.04%  "
$ws.Cells.Item(26, 4).Value = "'27.09"
$ws.Cells.Item(26, 5).Value = "  -1.69%  "
$ws.Cells.Item(27, 5).Value = "  -0.03%  "
$ws.Cells.Item(28, 5).Value = "  +12.88%  "
$ws.Cells.Item(29, 4).Value = "'10.17"
$ws.Cells.Item(29, 5).Value = "  +1.11%  "
$ws.Cells.Item(30, 4).Value = "'38.02"
$ws.Cells.Item(30, 5).Value = "  +0.45%  "
$ws.Cells.Item(31, 4).Value = "'5.92"
$ws.Cells.Item(31, 5).Value = "  -1.03%  "
$ws.Cells.Item(32, 4).Value = "'155.80"
$ws.Cells.Item(32, 5).Value = "  -0.10%  "
$ws.Cells.Item(33, 4).Value = "'19.43"
$ws.Cells.Item(33, 5).Value = "  +4.84%  "
$ws.Cells.Item(34, 4).Value = "'3.35"
$ws.Cells.Item(34, 5).Value = "  +1.43%  "
$ws.Cells.Item(35, 4).Value = "'2.09"
$ws.Cells.Item(35, 5).Value = "  -3.94%  "
$ws.Cells.Item(36, 5).Value = "  -2.11%  "
$ws.Cells.Item(37, 4).Value = "'2.61"
$ws.Cells.Item(37, 5).Value = "  -4.83%  "
$ws.Cells.Item(38, 5).Value = "  -0.60%  "
$ws.Cells.Item(39, 5).Value = "  +1.26%  "
$ws.Cells.Item(40, 4).Value = "'24.24"
$ws.Cells.Item(40, 5).Value = "  -8.09%  "
$ws.Cells.Item(41, 4).Value = "'3.40"
$ws.Cells.Item(41, 5).Value = "  +0.34%  "
$ws.Cells.Item(42, 5).Value = "  +0.28%  "
$ws.Cells.Item(43, 5).Value = "  -2.30%  "
$ws.Cells.Item(44, 4).Value = "'0.0304"
$ws.Cells.Item(44, 5).Value = "  +0.79%  "
$ws.Cells.Item(45, 5).Value = "  +0.11%  "
$ws.Cells.Item(46, 4).Value = "2.031.02"
$ws.Cells.Item(46, 5).Value = "  -1.74%  "
$ws.Cells.Item(47, 4).Value = "'84.63"
$ws.Cells.Item(47, 5).Value = "  -3.62%  "
$ws.Cells.Item(48, 4).Value = "'8.99"
$ws.Cells.Item(48, 5).Value = "  -2.28%  "
$ws.Cells.Item(49, 4).Value = "'74.94"
$ws.Cells.Item(49, 5).Value = "  -0.24%  "
$ws.Cells.Item(50, 4).Value = "2.763.19"
$ws.Cells.Item(50, 5).Value = "  -1.30%  "
$ws.Cells.Item(51, 5).Value = "  +0.84%  "
